$d = $word.ActiveDocument

$old = "Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός Λέων 2022: 14-23 Απριλίου, 14-23 Μαΐου"
$new = "2022 Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός Λέων: 14-23 Απριλίου, 14-23 Μαΐου"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
